# Updates odds/score-distribution figures scraped for the weekly FlashScore
# fixtures sheet (2025-05-05), plus one kickoff-time correction.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: BRAZIL - SERIE A BETANO | Bragantino vs Mirassol
$ws.Range("J9").Value = 1.07
$ws.Range("K9").Value = 9
$ws.Range("N9").Value = 2.08
$ws.Range("O9").Value = 1.73

# Row 10: BRAZIL - SERIE A BETANO | Juventude vs Atletico-MG
$ws.Range("G10").Value = 3.3
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 2.2
$ws.Range("U10").Value = 15
$ws.Range("Z10").Value = 9
$ws.Range("AG10").Value = 9.5

# Row 11: BRAZIL - SERIE C | Figueirense vs SER Caxias
$ws.Range("G11").Value = 2.02
$ws.Range("Q11").Value = 2.35

# Row 22: COLOMBIA - PRIMERA A | Chico vs Envigado
$ws.Range("G22").Value = 2.15
$ws.Range("I22").Value = 4.33
$ws.Range("AC22").Value = 67
$ws.Range("AF22").Value = 19
$ws.Range("AH22").Value = 41

# Row 26: COLOMBIA - PRIMERA B | Huila vs Patriotas
$ws.Range("H26").Value = 2.72
$ws.Range("I26").Value = 3.1
$ws.Range("J26").Value = 1.12
$ws.Range("K26").Value = 5.6
$ws.Range("L26").Value = 1.5
$ws.Range("M26").Value = 2.25
$ws.Range("N26").Value = 2.42
$ws.Range("O26").Value = 1.44
$ws.Range("P26").Value = 1.6
$ws.Range("Q26").Value = 2.07
$ws.Range("R26").Value = 1.98
$ws.Range("S26").Value = 1.65
$ws.Range("T26").Value = 6.3
$ws.Range("U26").Value = 11.25
$ws.Range("V26").Value = 9.75
$ws.Range("W26").Value = 28
$ws.Range("Y26").Value = 40
$ws.Range("Z26").Value = 6
$ws.Range("AA26").Value = 5.5
$ws.Range("AB26").Value = 16.5
$ws.Range("AC26").Value = 110
$ws.Range("AE26").Value = 7
$ws.Range("AF26").Value = 14.5
$ws.Range("AG26").Value = 11.75
$ws.Range("AH26").Value = 45
$ws.Range("AI26").Value = 35

# Row 27: CYPRUS - CYPRUS LEAGUE | Anorthosis vs Paralimni
$ws.Range("G27").Value = 3
$ws.Range("I27").Value = 2.15
$ws.Range("U27").Value = 16
$ws.Range("V27").Value = 10.75
$ws.Range("W27").Value = 37
$ws.Range("X27").Value = 25
$ws.Range("Y27").Value = 32
$ws.Range("Z27").Value = 10.75
$ws.Range("AB27").Value = 13.5
$ws.Range("AE27").Value = 8.25
$ws.Range("AF27").Value = 10.75
$ws.Range("AH27").Value = 20
$ws.Range("AI27").Value = 17

# Row 28: CYPRUS - CYPRUS LEAGUE | Karmiotissa vs Omonia Aradippou
$ws.Range("G28").Value = 1.72
$ws.Range("I28").Value = 4.75
$ws.Range("R28").Value = 1.7
$ws.Range("S28").Value = 1.93
$ws.Range("W28").Value = 14.5
$ws.Range("AB28").Value = 14
$ws.Range("AF28").Value = 28

# Row 29: CYPRUS - CYPRUS LEAGUE | Omonia 29th May vs Nea Salamis
$ws.Range("G29").Value = 3.6
$ws.Range("T29").Value = 13
$ws.Range("Y29").Value = 35
$ws.Range("AA29").Value = 7.2
$ws.Range("AB29").Value = 13

# Row 30: DENMARK - SUPERLIGA | Sonderjyske vs Aalborg
$ws.Range("I30").Value = 3.1
$ws.Range("J30").Value = 1.04
$ws.Range("K30").Value = 12
$ws.Range("X30").Value = 19
$ws.Range("AF30").Value = 15
$ws.Range("AH30").Value = 29

# Row 32: EGYPT - PREMIER LEAGUE | Petrojet vs Ceramica Cleopatra
$ws.Range("G32").Value = 3.7
$ws.Range("I32").Value = 2.15
$ws.Range("K32").Value = 6.2
$ws.Range("N32").Value = 2.1
$ws.Range("S32").Value = 1.91
$ws.Range("U32").Value = 20
$ws.Range("X32").Value = 35
$ws.Range("Z32").Value = 6.2
$ws.Range("AF32").Value = 10

# Row 40: INDONESIA - LIGA 1 | Persik Kediri vs Persebaya
$ws.Range("G40").Value = 3.05
$ws.Range("I40").Value = 2.2
$ws.Range("M40").Value = 2.87
$ws.Range("N40").Value = 1.93
$ws.Range("P40").Value = 1.43
$ws.Range("Q40").Value = 2.6
$ws.Range("R40").Value = 1.75
$ws.Range("S40").Value = 1.87
$ws.Range("T40").Value = 9
$ws.Range("U40").Value = 15.5
$ws.Range("V40").Value = 11
$ws.Range("W40").Value = 37
$ws.Range("X40").Value = 27
$ws.Range("Y40").Value = 35
$ws.Range("Z40").Value = 9.25
$ws.Range("AB40").Value = 14.5
$ws.Range("AC40").Value = 70
$ws.Range("AD40").Value = 600
$ws.Range("AE40").Value = 7.4
$ws.Range("AF40").Value = 10.5
$ws.Range("AG40").Value = 9
$ws.Range("AH40").Value = 21
$ws.Range("AI40").Value = 18.5

# Row 41: INDONESIA - LIGA 1 | Persis Solo vs Arema FC
$ws.Range("G41").Value = 1.7
$ws.Range("I41").Value = 4.35
$ws.Range("T41").Value = 8
$ws.Range("U41").Value = 8.75
$ws.Range("W41").Value = 14
$ws.Range("X41").Value = 13
$ws.Range("Y41").Value = 22
$ws.Range("Z41").Value = 12.5
$ws.Range("AB41").Value = 13.5
$ws.Range("AE41").Value = 14.5
$ws.Range("AF41").Value = 27
$ws.Range("AH41").Value = 75
$ws.Range("AI41").Value = 37
$ws.Range("AJ41").Value = 37

# Row 44: IRELAND - PREMIER DIVISION | Shelbourne vs Waterford
$ws.Range("J44").Value = 1.05
$ws.Range("L44").Value = 1.33

# Row 45: IRELAND - PREMIER DIVISION | Shamrock Rovers vs Sligo Rovers
$ws.Range("J45").Value = 1.03
$ws.Range("L45").Value = 1.19

# Row 46: IRELAND - PREMIER DIVISION | Derry City vs St. Patricks
$ws.Range("J46").Value = 1.05
$ws.Range("L46").Value = 1.37

# Row 47: ISRAEL - LIGAT HA'AL | Maccabi Tel Aviv vs H. Beer Sheva
$ws.Range("J47").Value = 1.03
$ws.Range("L47").Value = 1.17

# Row 51: PARAGUAY - DIVISION INTERMEDIA | Sol de America vs Resistencia
$ws.Range("O51").Value = 1.37

# Row 54: POLAND - DIVISION 1 | Kotwica Kolobrzeg vs Chrobry Glogow
$ws.Range("K54").Value = 10

# Row 55: PORTUGAL - LIGA PORTUGAL | Rio Ave vs Estrela
$ws.Range("C55").Value = "15:30"

# Row 60: ROMANIA - LIGA 2 | Csikszereda M. Ciuc vs Metaloglobus Bucharest
$ws.Range("G60").Value = 1.7
$ws.Range("H60").Value = 3.7
$ws.Range("I60").Value = 4.8
$ws.Range("L60").Value = 1.3
$ws.Range("M60").Value = 3.35
$ws.Range("N60").Value = 1.88
$ws.Range("O60").Value = 1.85
$ws.Range("R60").Value = 1.85
$ws.Range("S60").Value = 1.85
$ws.Range("T60").Value = 6.5
$ws.Range("U60").Value = 8
$ws.Range("V60").Value = 8.75
$ws.Range("W60").Value = 13.5
$ws.Range("X60").Value = 14.5
$ws.Range("Y60").Value = 30
$ws.Range("AA60").Value = 7.5
$ws.Range("AB60").Value = 17.5
$ws.Range("AC60").Value = 90
$ws.Range("AD60").Value = 800
$ws.Range("AE60").Value = 12.5
$ws.Range("AF60").Value = 30
$ws.Range("AG60").Value = 16.5
$ws.Range("AH60").Value = 100
$ws.Range("AI60").Value = 55
$ws.Range("AJ60").Value = 55

# Row 63: SAUDI ARABIA - DIVISION 1 | Al Batin vs Al Adalah
$ws.Range("G63").Value = 3.35
$ws.Range("H63").Value = 3.35
$ws.Range("O63").Value = 1.93
$ws.Range("R63").Value = 1.57
$ws.Range("S63").Value = 2.1
$ws.Range("T63").Value = 12
$ws.Range("U63").Value = 20
$ws.Range("V63").Value = 11.5
$ws.Range("W63").Value = 50
$ws.Range("AD63").Value = 250
$ws.Range("AE63").Value = 8.5

# Row 65: SOUTH KOREA - K LEAGUE 1 | Gwangju FC vs Gimcheon Sangmu
$ws.Range("N65").Value = 2.2
$ws.Range("O65").Value = 1.65

# Row 66: SOUTH KOREA - K LEAGUE 1 | Ulsan HD vs Pohang
$ws.Range("N66").Value = 1.93
$ws.Range("O66").Value = 1.93

# Row 68: SWEDEN - ALLSVENSKAN | Malmo FF vs Brommapojkarna
$ws.Range("G68").Value = 1.44
$ws.Range("K68").Value = 17
$ws.Range("N68").Value = 1.6
$ws.Range("O68").Value = 2.3
$ws.Range("W68").Value = 10
$ws.Range("AA68").Value = 9
$ws.Range("AH68").Value = 67

# Row 74: UNITED ARAB EMIRATES - UAE LEAGUE | Ajman vs Al Urooba
$ws.Range("G74").Value = 1.82
$ws.Range("H74").Value = 3.8
$ws.Range("I74").Value = 3.75
$ws.Range("K74").Value = 8.5
$ws.Range("P74").Value = 1.33
$ws.Range("Q74").Value = 3.05
$ws.Range("R74").Value = 1.65
$ws.Range("S74").Value = 2.12
$ws.Range("T74").Value = 8.5
$ws.Range("U74").Value = 9.5
$ws.Range("W74").Value = 15
$ws.Range("Z74").Value = 8.5
$ws.Range("AA74").Value = 7.5
$ws.Range("AB74").Value = 13.5
$ws.Range("AE74").Value = 13.5
